# tkinter_test.xlsx - rename PDFs by order date + company name
# Updates header fields, rebuilds the order-lines table (rows 5-15) with
# new data, and re-merges the Model(A)/Suma(B) columns per item group.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header block -----------------------------------------------------
$ws.Range("A1").Value = "FIRMA2"
$ws.Range("B2").Value = "13.5.2020"
$ws.Range("E2").ClearContents()
$ws.Range("B3").Value = "14.4.2020"
$ws.Range("E3").Value = 2567

# --- Style templates for the new table rows ----------------------------
# Row 5 already carries the "header row" formatting we need to replicate
# (A: s4, B: s5, C: s6, D: s5, E: s7); E5 carries the plain "continuation"
# style (s7) used by A/B on rows that are merged into the row above.
$ws.Range("A5:E5").Copy() | Out-Null
$ws.Range("A9:E9").PasteSpecial(-4122) | Out-Null
$ws.Range("A11:E11").PasteSpecial(-4122) | Out-Null
$ws.Range("A12:E12").PasteSpecial(-4122) | Out-Null
$ws.Range("A14:E14").PasteSpecial(-4122) | Out-Null

$ws.Range("C5:E5").Copy() | Out-Null
$ws.Range("C6:E6").PasteSpecial(-4122) | Out-Null
$ws.Range("C7:E7").PasteSpecial(-4122) | Out-Null
$ws.Range("C8:E8").PasteSpecial(-4122) | Out-Null
$ws.Range("C10:E10").PasteSpecial(-4122) | Out-Null
$ws.Range("C13:E13").PasteSpecial(-4122) | Out-Null
$ws.Range("C15:E15").PasteSpecial(-4122) | Out-Null

$ws.Range("E5").Copy() | Out-Null
$ws.Range("A6:B8").PasteSpecial(-4122) | Out-Null
$ws.Range("A10:B10").PasteSpecial(-4122) | Out-Null
$ws.Range("A13:B13").PasteSpecial(-4122) | Out-Null
$ws.Range("A15:B15").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# --- Table values --------------------------------------------------------
# Row 5 - D11 / styropian
$ws.Range("A5").Value = "D11"
$ws.Range("B5").Value = 241
$ws.Range("C5").Value = 12
$ws.Range("D5").Value = "styropian"

# Row 6 - ecru
$ws.Range("C6").Value = 34
$ws.Range("D6").Value = "ecru"

# Row 7 - czerwony
$ws.Range("C7").Value = 45
$ws.Range("D7").Value = "czerwony"

# Row 8 - czarny
$ws.Range("C8").Value = 150
$ws.Range("D8").Value = "czarny"

# Row 9 - B1 / styropian
$ws.Range("A9").Value = "B1"
$ws.Range("B9").Value = 73
$ws.Range("C9").Value = 33
$ws.Range("D9").Value = "styropian"

# Row 10 - czarny
$ws.Range("C10").Value = 40
$ws.Range("D10").Value = "czarny"

# Row 11 - M1 / styropian
$ws.Range("A11").Value = "M1"
$ws.Range("B11").Value = 50
$ws.Range("C11").Value = 50
$ws.Range("D11").Value = "styropian"

# Row 12 - Statyw metalowy / Slimak
$ws.Range("A12").Value = "Statyw metalowy"
$ws.Range("B12").Value = 110
$ws.Range("C12").Value = 60
$ws.Range("D12").Value = "Ślimak"

# Row 13 - 90
$ws.Range("C13").Value = 50
$ws.Range("D13").Value = "90"

# Row 14 - Statyw drewniany / bialy
$ws.Range("A14").Value = "Statyw drewniany"
$ws.Range("B14").Value = 70
$ws.Range("C14").Value = 20
$ws.Range("D14").Value = "biały"

# Row 15 - naturalny
$ws.Range("C15").Value = 50
$ws.Range("D15").Value = "naturalny"

# --- Re-merge the Model/Suma columns per item group -----------------------
$ws.Range("A5:A8").Merge() | Out-Null
$ws.Range("B5:B8").Merge() | Out-Null
$ws.Range("A9:A10").Merge() | Out-Null
$ws.Range("B9:B10").Merge() | Out-Null
$ws.Range("A11").Merge() | Out-Null
$ws.Range("B11").Merge() | Out-Null
$ws.Range("A12:A13").Merge() | Out-Null
$ws.Range("B12:B13").Merge() | Out-Null
$ws.Range("A14:A15").Merge() | Out-Null
$ws.Range("B14:B15").Merge() | Out-Null
